$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (plain decimals like "0.9998").
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = "26.899.87"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "1.741.44"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "310.77"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.5061"
$ws.Range("E7").Value = "  +8.09%  "
$ws.Range("D8").Value = "0.3572"
$ws.Range("E8").Value = "  +4.11%  "
$ws.Range("D9").Value = "42.02"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").Value = "0.07234"
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("D11").Value = "1.060"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").Value = "0.9996"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("E13").Value = "  +1.12%  "
$ws.Range("D14").Value = "5.951"
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("D15").Value = "1.741.96"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "6.802"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").Value = "86.59"
$ws.Range("E17").Value = "  -2.35%  "
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").Value = "0.06408"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "16.53"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "5.743"
$ws.Range("E22").Value = "  +1.91%  "
$ws.Range("D23").Value = "26.967.01"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("D24").Value = "11.23"
$ws.Range("E24").Value = "  +3.17%  "
$ws.Range("E25").Value = "  -4.45%  "
$ws.Range("D26").Value = "152.12"
$ws.Range("E26").Value = "  -3.03%  "
$ws.Range("E27").Value = "  +1.97%  "
$ws.Range("D28").Value = "1.941.20"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").Value = "2.209"
$ws.Range("E29").Value = "  +3.75%  "
$ws.Range("D30").Value = "119.30"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").Value = "1.039"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("D32").Value = "0.09579"
$ws.Range("E32").Value = "  +4.66%  "
$ws.Range("D33").Value = "3.578"
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("D34").Value = "5.361"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "0.05893"
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("D36").Value = "0.02177"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").Value = "10.97"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").Value = "0.1999"
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").Value = "1.422"
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.6035"
$ws.Range("E40").Value = "  +2.08%  "
$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").Value = "4.751"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").Value = "1.106"
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("D43").Value = "7.612"
$ws.Range("E43").Value = "  +1.93%  "
$ws.Range("D44").Value = "12.93"
$ws.Range("E44").Value = "  +2.57%  "
$ws.Range("D45").Value = "3.591"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").Value = "0.5650"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "119.92"
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("D48").Value = "1.843"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").Value = "0.9999"
$ws.Range("E51").Value = "  -0.05%  "
